$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new rows to make room ---
# Two new postings go in at the top (rows 2-3), pushing old rows 2-5 down to 4-7
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()
# One more new posting is inserted before the old row 4 (now sitting at row 6)
$ws.Rows.Item(6).Insert()

# --- Write cell values for rows 2-8 (row 1 header is unchanged) ---
# Row 2
$ws.Range("A2").Value = "2026-01-09 12:40:16"
$ws.Range("B2").Value = "製造業向け図面自動生成システムの開発・ツール化を支援してくださるエンジニア募集(AI/バックエンド)"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5460562"
$ws.Range("G2").Value = 435
$ws.Range("H2").Value = "🔥AI,Ai ◆ツール,開発"

# Row 3
$ws.Range("A3").Value = "2026-01-09 12:40:16"
$ws.Range("B3").Value = "施設管理・現場業務向け チェックリスト業務の自動化・報告書作成システム開発エンジニア募集"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5460563"
$ws.Range("G3").Value = 220
$ws.Range("H3").Value = "◆開発,システム開発 ◇管理"

# Row 4
$ws.Range("A4").Value = "2026-01-09 12:40:16"
$ws.Range("B4").Value = "初回 急募 自動カートインツール 開発のプロフェッショナルを探しています"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5467745"
$ws.Range("G4").Value = 120
$ws.Range("H4").Value = "◆ツール,開発"

# Row 5
$ws.Range("A5").Value = "2026-01-09 12:40:16"
$ws.Range("B5").Value = "【緊急募集】動画解析アプリ開発のプロフェッショナル"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5467910"
$ws.Range("G5").Value = 88
$ws.Range("H5").Value = "◆開発 ◇アプリ"

# Row 6
$ws.Range("A6").Value = "2026-01-09 12:40:16"
$ws.Range("B6").Value = "【急募】BtoB向け越境ECプラットフォーム開発のパートナー募集"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "3,000,000 円 ~ 5,000,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5468347"
$ws.Range("G6").Value = 75
$ws.Range("H6").Value = "◆開発"

# Row 7
$ws.Range("A7").Value = "2026-01-09 12:40:16"
$ws.Range("B7").Value = "【急募】大手保険システム会社でのPJ推進支援(PM・PL経験者募集/都内常駐)"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5467981"
$ws.Range("G7").Value = 40

# Row 8
$ws.Range("A8").Value = "2026-01-09 12:40:16"
$ws.Range("B8").Value = "限定公開 限定公開の仕事"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "5,000,000 円 ~ / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5467882"
$ws.Range("G8").Value = 25

# --- Rebuild hyperlinks on column F in row order so rIds come out sequential ---
$ws.Cells.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5460562")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5460563")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5467745")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5467910")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5468347")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5467981")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5467882")

# Hyperlinks.Add() creates a fresh duplicate "Hyperlink" cell style per call;
# re-apply the sheet's existing named style so every F-cell shares one style index
$ws.Range("F2:F8").Style = "Hyperlink"

# --- Column widths (xlsx <col> width = COM ColumnWidth + 0.8333333333333336) ---
$widthFudge = 0.8333333333333336
$ws.Columns.Item(2).ColumnWidth = 52 - $widthFudge   # B: 41 -> 52
$ws.Columns.Item(4).ColumnWidth = 32 - $widthFudge   # D: 30 -> 32
$ws.Columns.Item(8).ColumnWidth = 16 - $widthFudge   # H: 12 -> 16

